# Fix the data error in the item-collect table (rows 4-20, columns D:F)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("D4").Value = 4

# Row 5
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1

# Row 6
$ws.Range("E6").Value = 5

# Row 7
$ws.Range("D7").Value = 2
$ws.Range("F7").Value = 2

# Row 9
$ws.Range("E9").Value = 4

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1

# Row 11
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 2

# Row 12
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 2

# Row 13
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 3

# Row 14
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 3

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 4

# Row 16
$ws.Range("D16").Value = 5
$ws.Range("F16").Value = 2

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 2

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("F18").Value = 4

# Row 19
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 5

# Row 20
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 4

# Reset view: scroll to top and select D13 (matches the author's final cursor position)
$ws.Range("A1").Select()
$ws.Range("D13").Select()
